# Requisitos funcionales.xlsx - documentation update
#  - fix typo "genarales" -> "generales"
#  - bump version/date fields back to 0.1 / 2018-09-09 (43352)
#  - move the active selection to A10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Control interno" section header.
$ws.Range("A4").Value = "generales"

# Update version number and date fields at the top of the sheet.
$ws.Range("C1").Value = 0.1
$ws.Range("C2").Value = 43352

# Move the selection/active cell to A10, matching the saved view state.
$ws.Range("A10").Select()
